$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force-text helper pattern used inline for column D (numeric-looking strings)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "92.772.03"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.96%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.110.30"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.43%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.91"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.97%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "613.70"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.61%  "

$ws.Range("E7").Value = "  -2.23%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.390"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.30%  "

$ws.Range("E9").Value = "  -0.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.108.15"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.38%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.785"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +4.22%  "

$ws.Range("E12").Value = "  -3.64%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000245"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.95%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.508.02"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.93%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "33.94"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.95%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.43"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.34%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.116.63"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.33%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.81"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.41%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.53"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.88%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.84"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.68%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0000204"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "439.45"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.73%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.11"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.99%  "

$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.58"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -5.77%  "

$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "85.45"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -4.84%  "

$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.62"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.15%  "

$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.274.00"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.68%  "

$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.40%  "

$ws.Range("B30").Value = "Cronos"
$ws.Range("C30").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.182"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +8.09%  "

$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.126"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -10.73%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.234"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.40%  "

$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.03"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -30.41%  "

$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.15"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.22%  "

$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.09"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +7.69%  "

$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.163"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -5.20%  "

$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "25.70"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.61%  "

$ws.Range("B38").Value = "MantraDAO"
$ws.Range("C38").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.01"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +3.73%  "

$ws.Range("B39").Value = "PancakeSwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.90"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -8.68%  "

$ws.Range("B40").Value = "WhiteBITCoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.92"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +7.72%  "

$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.29"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.36%  "

$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "465.26"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.58%  "

$ws.Range("B43").Value = "PolygonEcosystemToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.428"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.54%  "

$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.34"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.93%  "

$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.05%  "

$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "159.28"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.83%  "

$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.682"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.44%  "

$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.83"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.99%  "

$ws.Range("B49").Value = "ImmutableX"
$ws.Range("C49").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.32"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.64%  "

$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "43.82"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.48%  "

$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0324"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.82%  "
